# Player.xlsx - "Property" sheet: modify queue lock, modify the way of saving player's data

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# --- 1. Unlock the queue: flip Save (column E) from TRUE to FALSE for rows 44-67 ---
for ($r = 44; $r -le 67; $r++) {
    $ws.Range("E" + $r).Value = $false
}

# --- 2. De-highlight rows 76 and 77 (they were marked with the red/yellow
#         "locked" style, now they should look like the regular rows, e.g. 78) ---

# Row 76: A76 should look like A68 (plain text-formatted style), the rest of
# the row should look like the corresponding cells on row 78 (no special style).
$ws.Range("A68").Copy()
$ws.Range("A76").PasteSpecial(-4122)

$ws.Range("B78").Copy()
$ws.Range("B76").PasteSpecial(-4122)

$ws.Range("G78").Copy()
$ws.Range("G76").PasteSpecial(-4122)

$ws.Range("H78").Copy()
$ws.Range("H76").PasteSpecial(-4122)

$ws.Range("I78").Copy()
$ws.Range("I76").PasteSpecial(-4122)

$ws.Range("J78").Copy()
$ws.Range("J76").PasteSpecial(-4122)

# Row 77: every formatted cell should look like the matching cell on row 78
# (no special style at all).
$ws.Range("A78").Copy()
$ws.Range("A77").PasteSpecial(-4122)

$ws.Range("B78").Copy()
$ws.Range("B77").PasteSpecial(-4122)

$ws.Range("G78").Copy()
$ws.Range("G77").PasteSpecial(-4122)

$ws.Range("H78").Copy()
$ws.Range("H77").PasteSpecial(-4122)

$ws.Range("I78").Copy()
$ws.Range("I77").PasteSpecial(-4122)

$ws.Range("J78").Copy()
$ws.Range("J77").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- 3. Move the sheet's active selection from E76 to H78 ---
$ws.Range("H78").Select() | Out-Null
